# "flujo con producto random y datos random"
#
# Adds a new "Banco" column to the Epayco sheet and a sample data row
# (mirroring the random sample data already present on the User sheet),
# then moves the active selection/tab to the Epayco sheet.

$wb = $excel.ActiveWorkbook
$wsUser   = $wb.Worksheets.Item("User")
$wsEpayco = $wb.Worksheets.Item("Epayco")

# --- Epayco sheet: new "Banco" header column ---
$wsEpayco.Range("F1").Value = "Banco"

# Copy the row-2 cell formatting (border + hyperlink style) from the User
# sheet's sample row so the new Epayco row matches the existing look.
$wsUser.Range("A2:E2").Copy()
$wsEpayco.Range("A2:E2").PasteSpecial(-4122) # xlPasteFormats

# --- Epayco sheet: new sample data row (random product / random data) ---
$wsEpayco.Range("A2").Value = "Juan"
$wsEpayco.Range("B2").Value = "uno"
$wsEpayco.Range("C2").Value = 888888881
$wsEpayco.Range("D2").Value = "Excel@excel.com"
$wsEpayco.Range("E2").Value = 3001111111
$wsEpayco.Range("F2").Value = 1001

# Hyperlink on the email cell, same as the equivalent cells on the User sheet.
[void]$wsEpayco.Hyperlinks.Add($wsEpayco.Range("D2"), "mailto:Excel@excel.com")

# Adding the hyperlink resets D2's style; reapply the bordered hyperlink
# style (matches D2:D7 on the User sheet).
$wsUser.Range("D2").Copy()
$wsEpayco.Range("D2").PasteSpecial(-4122) # xlPasteFormats

# --- Update selections ---
[void]$wsUser.Range("A2:E2").Select()
[void]$wsEpayco.Range("F4").Select()

# --- Epayco becomes the active/visible tab ---
$wsEpayco.Activate()
